$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# "Visitors": drop the 2020 column (column B), shifting the remaining
# year columns (2019/2018/2017) one to the left, and leave the sheet's
# selection on L18 (matches the state captured after this edit).
$wsVisitors = $wb.Worksheets.Item("Visitors")
$wsVisitors.Select()
$wsVisitors.Columns.Item(2).Delete()
$wsVisitors.Range("L18").Select()

# Remove the "NOTES" sheet entirely.
$wsNotes = $wb.Worksheets.Item("NOTES")
$wsNotes.Delete()

# Make "CollectionUseDelivery" the active/selected sheet (was "SpacesStaff").
$wb.Worksheets.Item("CollectionUseDelivery").Select()
